$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = 2; 3 = 0; 4 = 3; 5 = 0; 6 = 2; 7 = 0; 8 = 2; 9 = 1; 10 = 2;
    11 = 1; 12 = 2; 13 = 1; 14 = 0; 15 = 2; 16 = 0; 17 = 2; 18 = 0; 19 = 0;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 2; 26 = 2; 27 = 0; 28 = 1;
    29 = 1; 30 = 1; 31 = 0; 32 = 1; 33 = 0; 34 = 1; 35 = 1; 36 = 2; 37 = 3;
    38 = 2; 39 = 0; 40 = 1; 41 = 1; 42 = 1; 43 = 1; 44 = 1; 45 = 1; 46 = 1;
    47 = 2; 48 = 4; 49 = 2; 50 = 1; 51 = 2; 52 = 2; 53 = 1; 54 = 1; 55 = 3;
    56 = 2; 57 = 1; 58 = 2; 59 = 0; 60 = 4; 61 = 0; 62 = 3; 63 = 2; 64 = 1;
    65 = 0; 66 = 3; 67 = 1; 68 = 2; 69 = 0; 70 = 4; 71 = 0; 72 = 4; 73 = 2;
    74 = 1; 75 = 2; 76 = 2; 77 = 1; 78 = 2; 79 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
